$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = -12.99219999999999
$ws.Range("A9").Value = -21.89350000000002
$ws.Range("C12").Value = -11.4351
$ws.Range("A18").Value = -22.17060000000001
$ws.Range("A20").Value = -20.85489999999998
$ws.Range("C26").Value = -12.8153
$ws.Range("A27").Value = -21.952
$ws.Range("C27").Value = -13.1235
$ws.Range("C29").Value = -11.1741
$ws.Range("C37").Value = -14.49269999999998
$ws.Range("C38").Value = -12.9934
$ws.Range("C51").Value = -12.11989999999999
$ws.Range("C55").Value = -13.8018
$ws.Range("A69").Value = -21.9338
$ws.Range("C69").Value = -12.0327
$ws.Range("C70").Value = -11.4315
$ws.Range("A76").Value = -19.87489999999999
$ws.Range("A82").Value = -21.9766
$ws.Range("C83").Value = -13.77359999999999
$ws.Range("C102").Value = -13.4016
